# Updated cryptos list on Thu Oct 17 23:25:49 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data columns keep their original text formatting so that
# numeric-looking strings (e.g. "1.00", "591.86") are not silently
# reinterpreted as numbers by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.350.27"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.605.02"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "591.86"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").Value = "150.45"
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.544"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").Value = "2.603.13"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").Value = "0.128"
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("D13").Value = "0.343"
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").Value = "27.28"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").Value = "3.081.60"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "0.0000181"
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("D17").Value = "67.186.09"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "2.606.46"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "370.11"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "11.03"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("D22").Value = "4.28"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("E24").Value = "  -3.68%  "
$ws.Range("D25").Value = "73.14"
$ws.Range("E25").Value = "  +4.44%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "9.90"
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").Value = "2.735.47"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "577.09"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").Value = "0.0₃0984"
$ws.Range("E31").Value = "  -5.91%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "7.65"
$ws.Range("E32").Value = "  -3.52%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.35"
$ws.Range("E33").Value = "  -5.69%  "
$ws.Range("E34").Value = "  -3.08%  "
$ws.Range("E36").Value = "  -4.15%  "
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("D38").Value = "158.37"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").Value = "19.03"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").Value = "0.365"
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("D42").Value = "5.22"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("D44").Value = "17.12"
$ws.Range("E44").Value = "  +4.15%  "
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "153.13"
$ws.Range("E46").Value = "  -2.26%  "
$ws.Range("D47").Value = "0.0₆0282"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "3.70"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0777"
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "1.67"
$ws.Range("E50").Value = "  -4.69%  "
$ws.Range("D51").Value = "21.31"
$ws.Range("E51").Value = "  +1.36%  "
